$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Append new sedimentary rock rows (breccia..limestoneDesc) at the bottom first
# so their shared strings get indices 50-69 (matches diff ordering).
$ws.Range("A28").Value = 'breccia'
$ws.Range("B28").Value = 'Breccia'
$ws.Range("C28").Value = 0.6

$ws.Range("A29").Value = 'brecciaDesc'
$ws.Range("B29").Value = 'Sedimentary rock. Clastic'
$ws.Range("C29").Value = 5

$ws.Range("A30").Value = 'sandstone'
$ws.Range("B30").Value = 'Sandstone'
$ws.Range("C30").Value = 0.6

$ws.Range("A31").Value = 'sandstoneDesc'
$ws.Range("B31").Value = 'Sedimentary rock. Clastic'
$ws.Range("C31").Value = 5

$ws.Range("A32").Value = 'shale'
$ws.Range("B32").Value = 'Shale'
$ws.Range("C32").Value = 0.6

$ws.Range("A33").Value = 'shaleDesc'
$ws.Range("B33").Value = 'Sedimentary rock. Clastic'
$ws.Range("C33").Value = 5

$ws.Range("A34").Value = 'siltstone'
$ws.Range("B34").Value = 'Siltstone'
$ws.Range("C34").Value = 0.6

$ws.Range("A35").Value = 'siltstoneDesc'
$ws.Range("B35").Value = 'Sedimentary rock. Clastic'
$ws.Range("C35").Value = 5

$ws.Range("A36").Value = 'bituminousCoal'
$ws.Range("B36").Value = 'Bituminous Coal'
$ws.Range("C36").Value = 1

$ws.Range("A37").Value = 'bituminousCoalDesc'
$ws.Range("B37").Value = 'Sedimentary rock. Organic'
$ws.Range("C37").Value = 5

$ws.Range("A38").Value = 'limestone'
$ws.Range("B38").Value = 'Limestone'
$ws.Range("C38").Value = 0.5

$ws.Range("A39").Value = 'limestoneDesc'
$ws.Range("B39").Value = 'Sedimentary rock. Organic'
$ws.Range("C39").Value = 5

# Step 2: Insert 6 rows at row 8 for grain size entries (before olivine).
# Column A (keys) is filled top-to-bottom first, then column B (values),
# so shared strings get indices 70-75 (keys) then 76-81 (values).
$ws.Rows("8:13").Insert()
$ws.Range("A8").Value = 'grainSize_LargeVariant'
$ws.Range("A9").Value = 'grainSize_Sand'
$ws.Range("A10").Value = 'grainSize_Silt'
$ws.Range("A11").Value = 'grainSize_Clay'
$ws.Range("A12").Value = 'grainSize_FineCourseCrystal'
$ws.Range("A13").Value = 'grainSize_Coarse'
$ws.Range("B8").Value = 'Pebbles, cobbles, and boulders.'
$ws.Range("B9").Value = 'Sand'
$ws.Range("B10").Value = 'Silt'
$ws.Range("B11").Value = 'Clay'
$ws.Range("B12").Value = 'Fine to coarse crystals.'
$ws.Range("B13").Value = 'Microscopic to very coarse.'

# Step 3: Insert 4 rows at row 18 for calcite/carbon entries (after quartzDesc,
# before peridotite), so their shared strings get indices 82-89.
$ws.Rows("18:21").Insert()
$ws.Range("A18").Value = 'calcite'
$ws.Range("B18").Value = 'Calcite'
$ws.Range("C18").Value = 0.6

$ws.Range("A19").Value = 'calciteDesc'
$ws.Range("B19").Value = 'Shells, bones, etc.'
$ws.Range("C19").Value = 5

$ws.Range("A20").Value = 'carbon'
$ws.Range("B20").Value = 'Carbon'
$ws.Range("C20").Value = 0.6

$ws.Range("A21").Value = 'carbonDesc'
$ws.Range("B21").Value = 'Dead plants and poop.'
$ws.Range("C21").Value = 5

# Step 4: Update sheet view to match final selection/scroll position
$ws.Range("B22").Select()
$excel.ActiveWindow.ScrollRow = 7
